$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 304.72726
$ws.Cells.Item(2, 9).Value = 92.71429000000001
$ws.Cells.Item(2, 10).Value = 675.75
$ws.Cells.Item(2, 11).Value = 92.71429000000001
$ws.Cells.Item(2, 12).Value = 675.75
$ws.Cells.Item(2, 13).Value = 20.28570999999999
$ws.Cells.Item(2, 14).Value = -901.75
$ws.Cells.Item(64, 8).Value = 7649.9
$ws.Cells.Item(64, 9).Value = 4749
$ws.Cells.Item(64, 11).Value = 4749
$ws.Cells.Item(64, 13).Value = -4501
$ws.Cells.Item(67, 8).Value = 7649.9
$ws.Cells.Item(67, 9).Value = 4749
$ws.Cells.Item(67, 11).Value = 4749
$ws.Cells.Item(67, 13).Value = -3891
$ws.Cells.Item(74, 8).Value = 6608.8887
$ws.Cells.Item(74, 9).Value = 6608.8887
$ws.Cells.Item(74, 11).Value = 6608.8887
$ws.Cells.Item(74, 13).Value = -5672.8887
$ws.Cells.Item(77, 8).Value = 6608.8887
$ws.Cells.Item(77, 9).Value = 6608.8887
$ws.Cells.Item(77, 11).Value = 33044.4435
$ws.Cells.Item(77, 13).Value = -28364.4435
$ws.Cells.Item(80, 8).Value = 111132536
$ws.Cells.Item(80, 9).Value = 200001490
$ws.Cells.Item(80, 10).Value = 46352
$ws.Cells.Item(80, 11).Value = 600004470
$ws.Cells.Item(80, 12).Value = 139056
$ws.Cells.Item(80, 13).Value = -600003472
$ws.Cells.Item(80, 14).Value = -141052
$ws.Cells.Item(83, 8).Value = 111132536
$ws.Cells.Item(83, 9).Value = 200001490
$ws.Cells.Item(83, 10).Value = 46352
$ws.Cells.Item(83, 11).Value = 1800013410
$ws.Cells.Item(83, 12).Value = 417168
$ws.Cells.Item(83, 13).Value = -1800008418
$ws.Cells.Item(83, 14).Value = -427152
$ws.Cells.Item(86, 8).Value = 214287660
$ws.Cells.Item(86, 10).Value = 125001320
$ws.Cells.Item(86, 12).Value = 125001320
$ws.Cells.Item(86, 14).Value = -125003566
$ws.Cells.Item(89, 8).Value = 214287660
$ws.Cells.Item(89, 10).Value = 125001320
$ws.Cells.Item(89, 12).Value = 625006600
$ws.Cells.Item(89, 14).Value = -625017832
$ws.Cells.Item(92, 8).Value = 41667870
$ws.Cells.Item(92, 9).Value = 50001188
$ws.Cells.Item(92, 11).Value = 50001188
$ws.Cells.Item(92, 13).Value = -49999940
$ws.Cells.Item(103, 8).Value = 851.85187
$ws.Cells.Item(103, 9).Value = 463.4737
$ws.Cells.Item(103, 10).Value = 1774.25
$ws.Cells.Item(103, 11).Value = 1390.4211
$ws.Cells.Item(103, 12).Value = 5322.75
$ws.Cells.Item(103, 13).Value = -804.4211
$ws.Cells.Item(103, 14).Value = -6494.75
$ws.Cells.Item(132, 8).Value = 32312.334
$ws.Cells.Item(132, 9).Value = 37417.105
$ws.Cells.Item(132, 10).Value = 3725.6
$ws.Cells.Item(132, 11).Value = 112251.315
$ws.Cells.Item(132, 12).Value = 11176.8
$ws.Cells.Item(132, 13).Value = -109721.315
$ws.Cells.Item(132, 14).Value = -16236.8
$ws.Cells.Item(138, 8).Value = 2610.47
$ws.Cells.Item(138, 9).Value = 1317.5555
$ws.Cells.Item(138, 10).Value = 2738.3406
$ws.Cells.Item(138, 11).Value = 3952.6665
$ws.Cells.Item(138, 12).Value = 8215.0218
$ws.Cells.Item(138, 13).Value = 1187.3335
$ws.Cells.Item(138, 14).Value = -18495.0218

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 26963480
$ws.Cells.Item(32, 9).Value = 29824150
$ws.Cells.Item(32, 11).Value = 29824150
$ws.Cells.Item(32, 13).Value = -29823863
$ws.Cells.Item(61, 8).Value = 3902
$ws.Cells.Item(61, 9).Value = 3553.2856
$ws.Cells.Item(61, 10).Value = 4512.25
$ws.Cells.Item(61, 11).Value = 3553.2856
$ws.Cells.Item(61, 12).Value = 4512.25
$ws.Cells.Item(61, 13).Value = -3341.2856
$ws.Cells.Item(61, 14).Value = -4936.25
$ws.Cells.Item(74, 8).Value = 3512.353
$ws.Cells.Item(74, 9).Value = 3871.9
$ws.Cells.Item(74, 11).Value = 3871.9
$ws.Cells.Item(74, 13).Value = -2997.9
$ws.Cells.Item(77, 8).Value = 3512.353
$ws.Cells.Item(77, 9).Value = 3871.9
$ws.Cells.Item(77, 11).Value = 19359.5
$ws.Cells.Item(77, 13).Value = -14991.5
$ws.Cells.Item(97, 8).Value = 3505.0908
$ws.Cells.Item(97, 9).Value = 2430.625
$ws.Cells.Item(97, 11).Value = 2430.625
$ws.Cells.Item(97, 13).Value = -1934.625
$ws.Cells.Item(120, 8).Value = 62249.4
$ws.Cells.Item(120, 10).Value = 62249.4
$ws.Cells.Item(120, 12).Value = 62249.4
$ws.Cells.Item(120, 14).Value = -71925.39999999999
$ws.Cells.Item(136, 8).Value = 3902
$ws.Cells.Item(136, 9).Value = 3553.2856
$ws.Cells.Item(136, 10).Value = 4512.25
$ws.Cells.Item(136, 11).Value = 10659.8568
$ws.Cells.Item(136, 12).Value = 13536.75
$ws.Cells.Item(136, 13).Value = -8109.856800000001
$ws.Cells.Item(136, 14).Value = -18636.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(96, 8).Value = 100329.336
$ws.Cells.Item(96, 9).Value = 75000
$ws.Cells.Item(96, 10).Value = 112994
$ws.Cells.Item(96, 11).Value = 75000
$ws.Cells.Item(96, 12).Value = 112994
$ws.Cells.Item(96, 13).Value = -72254
$ws.Cells.Item(96, 14).Value = -118486
$ws.Cells.Item(105, 8).Value = 2045.9
$ws.Cells.Item(105, 9).Value = 1169
$ws.Cells.Item(105, 11).Value = 1169
$ws.Cells.Item(105, 13).Value = 578
$ws.Cells.Item(107, 8).Value = 3626.25
$ws.Cells.Item(107, 9).Value = 1020
$ws.Cells.Item(107, 10).Value = 3998.5715
$ws.Cells.Item(107, 11).Value = 1020
$ws.Cells.Item(107, 12).Value = 3998.5715
$ws.Cells.Item(107, 13).Value = 900
$ws.Cells.Item(107, 14).Value = -7838.5715
$ws.Cells.Item(117, 8).Value = 115022
$ws.Cells.Item(117, 10).Value = 115022
$ws.Cells.Item(117, 12).Value = 115022
$ws.Cells.Item(117, 14).Value = -124200
$ws.Cells.Item(122, 8).Value = 115664
$ws.Cells.Item(122, 10).Value = 115664
$ws.Cells.Item(122, 12).Value = 115664
$ws.Cells.Item(122, 14).Value = -125464

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 3496.8
$ws.Cells.Item(105, 9).Value = 3494.5
$ws.Cells.Item(105, 11).Value = 3494.5
$ws.Cells.Item(105, 13).Value = -1747.5
$ws.Cells.Item(141, 8).Value = 482576
$ws.Cells.Item(141, 10).Value = 482576
$ws.Cells.Item(141, 12).Value = 482576
$ws.Cells.Item(141, 14).Value = -492936

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 46.4375
$ws.Cells.Item(38, 9).Value = 35.18182
$ws.Cells.Item(38, 10).Value = 71.2
$ws.Cells.Item(38, 11).Value = 105.54546
$ws.Cells.Item(38, 12).Value = 213.6
$ws.Cells.Item(38, 13).Value = 241.45454
$ws.Cells.Item(38, 14).Value = -907.6
$ws.Cells.Item(55, 8).Value = 2666
$ws.Cells.Item(55, 9).Value = 1550
$ws.Cells.Item(55, 10).Value = 3503
$ws.Cells.Item(55, 11).Value = 4650
$ws.Cells.Item(55, 12).Value = 10509
$ws.Cells.Item(55, 13).Value = -4473
$ws.Cells.Item(55, 14).Value = -10863
$ws.Cells.Item(107, 8).Value = 910.70966
$ws.Cells.Item(107, 10).Value = 678.26666
$ws.Cells.Item(107, 12).Value = 2034.79998
$ws.Cells.Item(107, 14).Value = -5874.79998
$ws.Cells.Item(113, 8).Value = 1144.7778
$ws.Cells.Item(113, 9).Value = 663
$ws.Cells.Item(113, 10).Value = 4999
$ws.Cells.Item(113, 11).Value = 1989
$ws.Cells.Item(113, 12).Value = 14997
$ws.Cells.Item(113, 13).Value = 181
$ws.Cells.Item(113, 14).Value = -19337
$ws.Cells.Item(133, 8).Value = 3486.1667
$ws.Cells.Item(133, 9).Value = 2729.25
$ws.Cells.Item(133, 11).Value = 8187.75
$ws.Cells.Item(133, 13).Value = -3127.75
$ws.Cells.Item(140, 8).Value = 2992.8333
$ws.Cells.Item(140, 9).Value = 1989.25
$ws.Cells.Item(140, 11).Value = 5967.75
$ws.Cells.Item(140, 13).Value = -787.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 99332.664
$ws.Cells.Item(113, 10).Value = 99499
$ws.Cells.Item(113, 12).Value = 99499
$ws.Cells.Item(113, 14).Value = -103839
$ws.Cells.Item(124, 8).Value = 115193.75
$ws.Cells.Item(124, 10).Value = 115193.75
$ws.Cells.Item(124, 12).Value = 115193.75
$ws.Cells.Item(124, 14).Value = -125013.75
$ws.Cells.Item(126, 8).Value = 3971
$ws.Cells.Item(126, 10).Value = 3965.4
$ws.Cells.Item(126, 12).Value = 11896.2
$ws.Cells.Item(126, 14).Value = -16836.2
$ws.Cells.Item(132, 8).Value = 4004.0908
$ws.Cells.Item(132, 9).Value = 4004.5
$ws.Cells.Item(132, 11).Value = 12013.5
$ws.Cells.Item(132, 13).Value = -9483.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 9734
$ws.Cells.Item(7, 9).Value = 9651
$ws.Cells.Item(7, 10).Value = 9900
$ws.Cells.Item(7, 11).Value = 9651
$ws.Cells.Item(7, 12).Value = 9900
$ws.Cells.Item(7, 13).Value = -9539
$ws.Cells.Item(7, 14).Value = -10124
$ws.Cells.Item(40, 8).Value = 3638.8
$ws.Cells.Item(40, 9).Value = 2899.6667
$ws.Cells.Item(40, 11).Value = 2899.6667
$ws.Cells.Item(40, 13).Value = -2763.6667
$ws.Cells.Item(108, 8).Value = 36715.668
$ws.Cells.Item(108, 10).Value = 36715.668
$ws.Cells.Item(108, 12).Value = 36715.668
$ws.Cells.Item(108, 14).Value = -44395.668
$ws.Cells.Item(122, 8).Value = 30189.2
$ws.Cells.Item(122, 9).Value = 30189.2
$ws.Cells.Item(122, 11).Value = 90567.60000000001
$ws.Cells.Item(122, 13).Value = -88117.60000000001
$ws.Cells.Item(126, 8).Value = 9734
$ws.Cells.Item(126, 9).Value = 9651
$ws.Cells.Item(126, 10).Value = 9900
$ws.Cells.Item(126, 11).Value = 28953
$ws.Cells.Item(126, 12).Value = 29700
$ws.Cells.Item(126, 13).Value = -26483
$ws.Cells.Item(126, 14).Value = -34640
$ws.Cells.Item(132, 8).Value = 4877.5483
$ws.Cells.Item(132, 9).Value = 4485.6313
$ws.Cells.Item(132, 10).Value = 5498.0835
$ws.Cells.Item(132, 11).Value = 13456.8939
$ws.Cells.Item(132, 12).Value = 16494.2505
$ws.Cells.Item(132, 13).Value = -10926.8939
$ws.Cells.Item(132, 14).Value = -21554.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 6500
$ws.Cells.Item(126, 9).Value = 6500
$ws.Cells.Item(126, 11).Value = 19500
$ws.Cells.Item(126, 13).Value = -17030
$ws.Cells.Item(136, 8).Value = 16751868
$ws.Cells.Item(136, 9).Value = 2076.625
$ws.Cells.Item(136, 11).Value = 6229.875
$ws.Cells.Item(136, 13).Value = -3679.875
